# Update column C ("Förändrad") date values from 2023-10-05 (45204)
# to 2023-10-08 (45207) for all data rows (2 through 33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 33
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45207
    }
}
